$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 771
$ws1.Range("F5").Value = 2500
$ws1.Range("F7").Value = 1850
$ws1.Range("F8").Value = 3183
$ws1.Range("F9").Value = 198
$ws1.Range("F10").Value = 4676
$ws1.Range("F11").Value = 434
$ws1.Range("F12").Value = 258
$ws1.Range("F14").Value = 602
$ws1.Range("F16").Value = 10
$ws1.Range("F18").Value = 634
$ws1.Range("F19").Value = 277
$ws1.Range("F20").Value = 16
$ws1.Range("F23").Value = 325
$ws1.Range("F24").Value = 4652
$ws1.Range("F28").Value = 5560
$ws1.Range("F30").Value = 1169
$ws1.Range("F31").Value = 231
$ws1.Range("F32").Value = 642
$ws1.Range("F33").Value = 4400
$ws1.Range("F37").Value = 765
$ws1.Range("F39").Value = 701
$ws1.Range("F40").Value = 706

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1070

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1070
$ws4.Range("F7").Value = 771
$ws4.Range("F8").Value = 2500
$ws4.Range("F10").Value = 1850
$ws4.Range("F12").Value = 3183
$ws4.Range("F13").Value = 198
$ws4.Range("F14").Value = 4676
$ws4.Range("F15").Value = 434
$ws4.Range("F16").Value = 258
$ws4.Range("F18").Value = 602
$ws4.Range("F20").Value = 10
$ws4.Range("F22").Value = 634
$ws4.Range("F23").Value = 277
$ws4.Range("F24").Value = 16
$ws4.Range("F28").Value = 325
$ws4.Range("F29").Value = 4652
$ws4.Range("F33").Value = 5560
$ws4.Range("F35").Value = 1169
$ws4.Range("F36").Value = 231
$ws4.Range("F37").Value = 642
$ws4.Range("F38").Value = 4400
$ws4.Range("F43").Value = 765
$ws4.Range("F45").Value = 701
$ws4.Range("F46").Value = 706

Write-Output "Updated F column interest counts across sheets."